$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ランサーズ")

# Update the "取得日時" (acquired timestamp) column for the newly appended
# batch of rows (2-8) to the new run's timestamp.
$newTimestamp = "2025-12-28 02:08:08"

for ($row = 2; $row -le 8; $row++) {
    $ws.Cells.Item($row, 1).Value = $newTimestamp
}
